$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.061.27'
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").Value = '3.859.89'
$ws.Range("E3").Value = '  +1.14%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '693.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.96'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.81%  '

$ws.Range("D7").Value = '3.856.22'
$ws.Range("E7").Value = '  +1.08%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.43'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +6.28%  '

$ws.Range("E12").Value = '  -0.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000259'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +5.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.68'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.89%  '

$ws.Range("D15").Value = '4.514.95'
$ws.Range("E15").Value = '  +1.17%  '

$ws.Range("D16").Value = '3.864.54'

$ws.Range("D17").Value = '71.155.82'
$ws.Range("E17").Value = '  +0.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.83'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.25'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.61%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.19'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '487.83'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.49%  '

$ws.Range("E23").Value = '  +0.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.74'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.75%  '

$ws.Range("E25").Value = '  +1.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.41'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.52'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.22%  '

$ws.Range("E28").Value = '  +0.72%  '

$ws.Range("D29").Value = '4.016.14'
$ws.Range("E29").Value = '  +1.18%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.12'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +8.50%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.62'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.30'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.79'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.181'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.28'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.92%  '

$ws.Range("D37").Value = '3.814.98'
$ws.Range("E37").Value = '  +1.09%  '

$ws.Range("E38").Value = '  +0.11%  '

$ws.Range("E39").Value = '  +1.43%  '

$ws.Range("E40").Value = '  +13.38%  '

$ws.Range("E41").Value = '  +0.01%  '

$ws.Range("E42").Value = '  +0.90%  '

$ws.Range("E43").Value = '  +4.70%  '

$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '164.27'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +4.36%  '

$ws.Range("E47").Value = '  +6.63%  '

$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.70'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.34%  '

$ws.Range("B49").Value = 'Arweave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.63'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.303'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.70'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.07%  '

